# Update column G ("K") values in Sheet1 to reflect the regenerated
# save_data using K instead of Strike#.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 3
    6  = 4
    7  = 1
    8  = 3
    9  = 1
    10 = 3
    11 = 1
    12 = 3
    13 = 5
    14 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
